# Reorders the player roster table (rows 2-18) on the active sheet so that
# the newly-available "Out of PO" free agents appear first, followed by the
# players that were already listed, keeping Cam Thomas / Donte DiVincenzo at
# the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing rows first so the shared-string table is rebuilt in the
# new write order (matches how the workbook was re-saved upstream).
$ws.Range("A2:C18").ClearContents()

$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Kel'el Ware", "C", "Miami Heat"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

# Write column-by-column (all player names, then all positions, then all
# teams) so the shared-string pool is rebuilt in the same first-seen order
# as the target workbook.
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $data[$i][2]
}
